$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, copy the date-format style from A52 to the new A53 cell before writing values,
# so the new row inherits the same number formatting as the rest of column A.
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Target data for rows 2..53: columns A (serial date), B (y_0 year), C (y_0_forecast),
# D (y_1 year), E (y_1_forecast). This reflects the corrected/regenerated forecast values,
# including a newly inserted leading observation and one additional trailing observation.
$data = @(
    @(2, 39400, 2007, 7.226520411029047, 2008, 8.987952903094421),
    @(3, 39583, 2008, 6.109216616889168, 2009, 8.045645122021906),
    @(4, 39765, 2008, 4.268860212333636, 2009, 6.77211531652997),
    @(5, 39948, 2009, -5.232639093663815, 2010, -0.9943400519801915),
    @(6, 40130, 2009, -7.266312015249776, 2010, -1.420242831007679),
    @(7, 40310, 2010, 2.12454018480297, 2011, 1.203634802640963),
    @(8, 40494, 2010, 6.958243460951929, 2011, 8.31992099499319),
    @(9, 40676, 2011, 8.081020954067775, 2012, 2.257871268432821),
    @(10, 40862, 2011, 9.469137444079934, 2012, 8.571528775834981),
    @(11, 41044, 2012, 4.489210662380949, 2013, 7.472658273721078),
    @(12, 41228, 2012, 3.358206407534947, 2013, 5.745831525574441),
    @(13, 41409, 2013, -0.8752093743685352, 2014, 1.998870338019265),
    @(14, 41592, 2013, 0.3081076735359067, 2014, 2.847379875994704),
    @(15, 41774, 2014, 4.891728508251214, 2015, 2.926340920335191),
    @(16, 41957, 2014, 3.901355411819707, 2015, 3.690459963535009),
    @(17, 42137, 2015, 4.818339085077583, 2016, 4.241902819910548),
    @(18, 42321, 2015, 5.331683351557981, 2016, 4.636575318346536),
    @(19, 42503, 2016, 4.067959312311897, 2017, 4.36243732366437),
    @(20, 42689, 2016, 3.254758369308375, 2017, 2.76788332063731),
    @(21, 42867, 2017, 4.613634856640747, 2018, 4.246555741688218),
    @(22, 43053, 2017, 5.246209615995667, 2018, 4.659862065074982),
    @(23, 43145, 2018, 6.155351106582874, 2019, 5.081201994458495),
    @(24, 43235, 2018, 3.898744563937395, 2019, 3.529300656691237),
    @(25, 43326, 2018, 5.021833280292598, 2019, 4.83330568333058),
    @(26, 43418, 2018, 4.86255966374296, 2019, 4.636196713604357),
    @(27, 43510, 2019, 4.237941638782527, 2020, 4.932073907517265),
    @(28, 43600, 2019, 3.841510956591465, 2020, 4.83848589746565),
    @(29, 43691, 2019, 3.042742667481324, 2020, 3.228901920070748),
    @(30, 43783, 2019, 2.764740011159428, 2020, 2.471557257221946),
    @(31, 43875, 2020, 2.059284312217757, 2021, 2.845541644111571),
    @(32, 43966, 2020, -0.8225206269755425, 2021, 1.159653508089242),
    @(33, 44068, 2020, -7.578477024949737, 2021, -5.048497756254311),
    @(34, 44159, 2020, -7.260793671746435, 2021, -1.387795042833839),
    @(35, 44251, 2021, -0.4318290737559183, 2022, -1.967992356539539),
    @(36, 44341, 2021, 4.409066926520455, 2022, 3.034999751677669),
    @(37, 44432, 2021, 4.379227219808146, 2022, 2.570052815033752),
    @(38, 44525, 2021, 4.097586525396268, 2022, 2.450242954096926),
    @(39, 44617, 2022, 6.557154773664364, 2023, 1.465110550383386),
    @(40, 44706, 2022, 8.064077385547574, 2023, 3.690055931494096),
    @(41, 44798, 2022, 7.397318165265498, 2023, 3.192017772210276),
    @(42, 44890, 2022, 7.824284864703746, 2023, 4.834990656989402),
    @(43, 44981, 2023, 1.287151040638124, 2024, 5.402455143891816),
    @(44, 45071, 2023, 0.3551698673347259, 2024, 4.038118345571751),
    @(45, 45163, 2023, -0.3046246622258053, 2024, 2.185848087053199),
    @(46, 45254, 2023, -1.24502235313334, 2024, -2.488220481262082),
    @(47, 45345, 2024, -3.197915208378399, 2025, 0.5664860188349996),
    @(48, 45436, 2024, -2.267078452724969, 2025, 0.5378929214800987),
    @(49, 45534, 2024, -2.567041707495976, 2025, -0.3140572462435154),
    @(50, 45618, 2024, -1.735114423676209, 2025, 1.056286187957367),
    @(51, 45713, 2025, 1.602010908728335, 2026, -1.246071640539481),
    @(52, 45800, 2025, 1.326993065386817, 2026, -0.2087957186147071),
    @(53, 45891, 2025, 2.64031107104763, 2026, 1.890175284700679)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

Write-Output "Applied forecast data update for rows 2-53."
